$wb = $excel.ActiveWorkbook

# Add a new worksheet for the news keywords and place it after "Лист1" (last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "news_keywords"

$values = @(
    "keyword",
    "News",
    "news",
    "новости",
    "novosti",
    "novosty"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $newSheet.Cells.Item($row, 1).Value = $values[$i]
}
